$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$ws2.Range("E3").Value = "2016-03-17 12:35:32"
$ws2.Range("H3").Value = "2016-03-17 12:35:50"

$ws3.Range("E3").Value = "2016-03-17 12:35:36"
$ws3.Range("H3").Value = "2016-03-17 12:35:56"
